# Flask-blogly part1: turn the Sheet2 "Tickets" mini-ERD table into a
# "Flights" table, add a "gender" column to Passengers, rename the
# Passengers PK from "id" to "passenger_id" and the Tickets/Flights PK
# from "id" to "flight_id", and pick up a couple of header fill colors
# that the other tables on Sheet2 were missing.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------
# 1. Passengers table (A:B) - rename PK, add gender row
# ---------------------------------------------------------------
$ws2.Range("B2").Value = "passenger_id"
$ws2.Range("B5").Value = "gender"

# ---------------------------------------------------------------
# 2. Rename the "Tickets" header to "Flights" and rebuild that
#    mini table as the Flights entity (PK flight_id, FK passenger_id,
#    FK airline, aircraft, flight_num, seat, departure, arrival,
#    FK from_city/from_country/to_city/to_country).
# ---------------------------------------------------------------
$ws2.Range("J1").Value = "Flights"
$ws2.Range("K2").Value = "flight_id"

$ws2.Range("J3").Value = "FK"
$ws2.Range("K3").Value = "passenger_id"

$ws2.Range("J4").Value = "FK"
$ws2.Range("K4").Value = "airline"

$ws2.Range("H4").Value = "hub_city"
$ws2.Range("K5").Value = "aircraft"

# Insert a row so "flight_num" gets its own slot and the remaining
# ticket fields (seat/departure/arrival/from_*/to_*) shift down one.
$ws2.Rows("6:6").Insert()
$ws2.Range("K6").Value = "flight_num"

# The from_*/to_* rows now live at 10-13 (shifted down by the insert
# above); just stamp the FK marker in column J for each.
$ws2.Range("J10").Value = "FK"
$ws2.Range("J11").Value = "FK"
$ws2.Range("J12").Value = "FK"
$ws2.Range("J13").Value = "FK"

# ---------------------------------------------------------------
# 3. Header fills: Passengers/Airports pick up fill colors matching
#    the Planets/Galaxies/Moons headers on Sheet1; Flights reuses the
#    Orbits sub-header color; Airline becomes yellow like the PK/FK
#    legend elsewhere in the workbook.
# ---------------------------------------------------------------
$ws1.Range("F7").Copy()
$ws2.Range("J1:K1").PasteSpecial(-4122)

$ws1.Range("I1").Copy()
$ws2.Range("J3").PasteSpecial(-4122)

$ws1.Range("B1").Copy()
$ws2.Range("G1:H1").PasteSpecial(-4122)
$ws2.Range("J4").PasteSpecial(-4122)

$ws1.Range("E3").Copy()
$ws2.Range("J2").PasteSpecial(-4122)
$ws2.Range("A6").PasteSpecial(-4122)

$ws2.Range("A1:B1").Interior.Color = 15652797
$ws2.Range("D1:E1").Interior.Color = 5296274

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 4. Selection / view state
# ---------------------------------------------------------------
$ws2.Activate()
$ws2.Range("D13").Select()
